$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Courtney", "Holcomb", "Courtney.Holcomb.1992@gmail.com", "J90PAVtH85i1vC2h", "81.28.96.172:15878", "a4ZJCKXpyPxqtgAt", "WZRYymEeVjfQeR6Z"),
    @("Jill", "Middleton", "Jill.Middleton.1980@gmail.com", "QkserjiJNZGdu5WF", "81.28.96.172:65448", "a4ZJCKXpyPxqtgAt", "WZRYymEeVjfQeR6Z"),
    @("Justine", "Lester", "Justine.Lester.1981@gmail.com", "DpzEuqQtJsn0l65F", "81.28.96.172:5596", "a4ZJCKXpyPxqtgAt", "WZRYymEeVjfQeR6Z"),
    @("Ola", "Hampton", "Ola.Hampton.1982@gmail.com", "kgQdKfc8IHlw1CWO", "81.28.96.172:65014", "a4ZJCKXpyPxqtgAt", "WZRYymEeVjfQeR6Z"),
    @("Earlene", "Simon", "Earlene.Simon.1983@gmail.com", "qlrLPymGRdcEZo1h", "81.28.96.172:3222", "a4ZJCKXpyPxqtgAt", "WZRYymEeVjfQeR6Z")
)

$startRow = 11
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
}

$ws.Rows.Item(11).Select()
